# The edit reorders the 5 data rows (rows 2-6) of the "Artfynd" sheet.
# Using full-row copy/paste (values only) via a scratch area so that
# text values that look like dates (e.g. "2023-08-15") are not
# reinterpreted/auto-converted by Excel, and so sparse (missing) cells
# correctly blank out destination cells that must become empty.
#
# Mapping of final row -> original row that its data comes from:
#   new row 2 <- old row 6
#   new row 3 <- old row 2
#   new row 4 <- old row 5
#   new row 5 <- old row 4
#   new row 6 <- old row 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues
$xlPasteValues = -4163

# Stash the current rows 2-6 into a scratch area (rows 100-104) so that
# the source data survives while we overwrite rows 2-6.
$ws.Range("A2:AY6").Copy()
$ws.Range("A100:AY104").PasteSpecial($xlPasteValues)

# Clear the destination rows completely first so columns that are blank
# in the source (no cell at all) truly end up blank in the destination
# instead of retaining old leftover content.
$ws.Range("A2:AY6").Clear()

# scratch row 100 = original row 2
# scratch row 101 = original row 3
# scratch row 102 = original row 4
# scratch row 103 = original row 5
# scratch row 104 = original row 6

$ws.Range("A104:AY104").Copy()
$ws.Range("A2:AY2").PasteSpecial($xlPasteValues)

$ws.Range("A100:AY100").Copy()
$ws.Range("A3:AY3").PasteSpecial($xlPasteValues)

$ws.Range("A103:AY103").Copy()
$ws.Range("A4:AY4").PasteSpecial($xlPasteValues)

$ws.Range("A102:AY102").Copy()
$ws.Range("A5:AY5").PasteSpecial($xlPasteValues)

$ws.Range("A101:AY101").Copy()
$ws.Range("A6:AY6").PasteSpecial($xlPasteValues)

# Clean up the scratch area.
$ws.Range("A100:AY104").Clear()
